{"js": "// Auto-generated: replace each unique cell/date text with its updated value.\nconst replacements = [\n  [\"2026-02-11 Wednesday\", \"2026-02-12 Thursday\"],\n  [\"75\u00d777=5775\", \"97\u00d730=2910\"],\n  [\"68\u00d735=2380\", \"23\u00d789=2047\"],\n  [\"69\u00d726=1794\", \"39\u00d735=1365\"],\n  [\"61\u00d753=3233\", \"17\u00d728=476\"],\n  [\"32\u00d757=1824\", \"72\u00d758=4176\"],\n  [\"90\u00d746=4140\", \"96\u00d791=8736\"],\n  [\"34\u00d762=2108\", \"83\u00d791=7553\"],\n  [\"45\u00d718=810\", \"48\u00d740=1920\"],\n  [\"95\u00d740=3800\", \"41\u00d783=3403\"],\n  [\"78\u00d726=2028\", \"88\u00d723=2024\"],\n  [\"72\u00d756=4032\", \"56\u00d770=3920\"],\n  [\"51\u00d718=918\", \"70\u00d777=5390\"],\n  [\"54\u00d775=4050\", \"20\u00d746=920\"],\n  [\"38\u00d749=1862\", \"21\u00d732=672\"],\n  [\"75\u00d719=1425\", \"96\u00d769=6624\"],\n  [\"95\u00d716=1520\", \"31\u00d780=2480\"],\n  [\"83\u00d766=5478\", \"74\u00d733=2442\"],\n  [\"79\u00d779=6241\", \"30\u00d721=630\"],\n  [\"38\u00d745=1710\", \"92\u00d736=3312\"],\n  [\"97\u00d748=4656\", \"60\u00d780=4800\"],\n  [\"68\u00d721=1428\", \"79\u00d737=2923\"],\n  [\"88\u00d793=8184\", \"71\u00d765=4615\"],\n  [\"94\u00d711=1034\", \"72\u00d758=4176\"],\n  [\"88\u00d717=1496\", \"66\u00d752=3432\"],\n  [\"38\u00d777=2926\", \"52\u00d769=3588\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2026-02-11 Wednesday', '2026-02-12 Thursday'),\n    @('75\u00d777=5775', '97\u00d730=2910'),\n    @('68\u00d735=2380', '23\u00d789=2047'),\n    @('69\u00d726=1794', '39\u00d735=1365'),\n    @('61\u00d753=3233', '17\u00d728=476'),\n    @('32\u00d757=1824', '72\u00d758=4176'),\n    @('90\u00d746=4140', '96\u00d791=8736'),\n    @('34\u00d762=2108', '83\u00d791=7553'),\n    @('45\u00d718=810', '48\u00d740=1920'),\n    @('95\u00d740=3800', '41\u00d783=3403'),\n    @('78\u00d726=2028', '88\u00d723=2024'),\n    @('72\u00d756=4032', '56\u00d770=3920'),\n    @('51\u00d718=918', '70\u00d777=5390'),\n    @('54\u00d775=4050', '20\u00d746=920'),\n    @('38\u00d749=1862', '21\u00d732=672'),\n    @('75\u00d719=1425', '96\u00d769=6624'),\n    @('95\u00d716=1520', '31\u00d780=2480'),\n    @('83\u00d766=5478', '74\u00d733=2442'),\n    @('79\u00d779=6241', '30\u00d721=630'),\n    @('38\u00d745=1710', '92\u00d736=3312'),\n    @('97\u00d748=4656', '60\u00d780=4800'),\n    @('68\u00d721=1428', '79\u00d737=2923'),\n    @('88\u00d793=8184', '71\u00d765=4615'),\n    @('94\u00d711=1034', '72\u00d758=4176'),\n    @('88\u00d717=1496', '66\u00d752=3432'),\n    @('38\u00d777=2926', '52\u00d769=3588'),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $rng = $d.Content\n    $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}"}
